$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.149.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.29%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.669.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.85%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.43%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'209.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.35%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.5239"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.03%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.45%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2625"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.86%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06343"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.66%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07538"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.75%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.677.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.64%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.446"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.51%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.5509"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.52%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'66.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.33%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.000007968"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -4.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.161.22"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.39%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -0.44%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.756"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.42%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'186.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.73%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.78%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.188"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.60%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'149.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.47%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -1.07%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.511"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.70%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'15.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.87%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.06415"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.64%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.350"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.57%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.276"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.33%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.514"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.64%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.414"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.16%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.645"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.46%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.74%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.6025"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.98%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.408"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.48%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.740"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.68%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.114.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.99%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'6.146"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.57%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01617"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.26%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.8666"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.72%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'100.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.14%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.823.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.56%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.72%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'55.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.40%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.26%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'8.077"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.10%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.05232"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.85%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.4239"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.14%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'5.930"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.54%  "
$ws.Range("E51").Style = "Normal"
